# Correção nos dados e inicio da analise PNAD 2009
# The B2 and F2 header cells held placeholder pandas-exported labels
# ("unnamed: 1_level_1" / "unnamed: 5_level_1"); fix them to "total",
# matching the other "total" sub-header already present in C2. This
# also removes the two now-unused shared strings from the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
